# Generate Report for Handoff
# Replace old localization UUID/hash with the new ones across all sheets,
# refresh the handoff timestamps, and clear the stale handback info.

$wb = $excel.ActiveWorkbook

$oldUuid = "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc"
$newUuid = "79e7da14-080e-4792-9dda-84539cd54f49"
$oldHash = "4ab74e330ad8179519427b2cff08588ae293166e"
$newHash = "802f4670e16e2bdd5f0faaa0f944c6ab45378cff"

# Same hyperlink target on every sheet - only the visible text changes.
$mdHyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/$oldUuid.md"

### ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid.md"

# B2 carries a hyperlink; drop the old one and recreate it so only the
# displayed text changes while the link target stays the same.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $mdHyperlinkTarget, [Type]::Missing, [Type]::Missing, "e2e\$newUuid.md") | Out-Null

$wsOverview.Range("G2").Value = "2016-08-30 11:08:39"

### ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdHyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newUuid.md") | Out-Null

$wsZhCn.Range("G2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 11:08:32"

# Latest Target File / Latest Handback File are no longer available.
$wsZhCn.Range("I2").ClearContents()
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").ClearContents()

$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Columns.Item(9).ColumnWidth = 17.833333333333336
$wsZhCn.Columns.Item(10).ColumnWidth = 20.833333333333336

### ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdHyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newUuid.md") | Out-Null

$wsDeDe.Range("G2").Value = "$newUuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-30 11:08:39"

$wsDeDe.Range("I2").ClearContents()
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").ClearContents()

$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Columns.Item(9).ColumnWidth = 17.833333333333336
$wsDeDe.Columns.Item(10).ColumnWidth = 20.833333333333336
